$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new 2022-Q3 data row at the
#    top of the data (row 2), shifting the existing quarters down by one row,
#    and append a new trailing row for what used to be the last row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Shift the B/C/D (quarter label, count, value) columns down one row each,
# working from the bottom up so we never clobber a value before reading it.
$summary.Range("B6").Value = $summary.Range("B5").Value()
$summary.Range("C6").Value = $summary.Range("C5").Value()
$summary.Range("D6").Value = $summary.Range("D5").Value()

$summary.Range("B5").Value = $summary.Range("B4").Value()
$summary.Range("C5").Value = $summary.Range("C4").Value()
$summary.Range("D5").Value = $summary.Range("D4").Value()

$summary.Range("B4").Value = $summary.Range("B3").Value()
$summary.Range("C4").Value = $summary.Range("C3").Value()
$summary.Range("D4").Value = $summary.Range("D3").Value()

$summary.Range("B3").Value = $summary.Range("B2").Value()
$summary.Range("C3").Value = $summary.Range("C2").Value()
$summary.Range("D3").Value = $summary.Range("D2").Value()

# New first data row: 2022-Q3
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.02

# Column A is a simple running index (0,1,2,3,...). It keeps its original
# values for the existing rows; only the new trailing row needs a value,
# copied (for style) from the row above then overwritten.
$summary.Range("A5").Copy($summary.Range("A6"))
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# 2. Insert a brand-new worksheet named "2022-Q3" right after "总计" holding
#    the fund detail rows for the new quarter. Cloning the sheet that is
#    currently in that slot ("2021-Q4") keeps all of its formatting
#    (sheetPr, header/index cell styles, page margins, ...) and places the
#    clone directly before it in tab order.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item(2)
$template.Copy($template)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Fund code / name / scale / position columns are stored as text in the
# source workbook (even the numeric-looking ones), so force text with a
# leading apostrophe to avoid Excel auto-converting them to numbers (which
# would also strip the leading zeros off the fund codes).
$newSheet.Range("B2").Value = "'000059"
$newSheet.Range("C2").Value = "国联安中证医药100指数A"
$newSheet.Range("D2").Value = "'1.70"
$newSheet.Range("E2").Value = "'92.19"
$newSheet.Range("F2").Value = "'1.17"
$newSheet.Range("G2").Value = "'0.0199"
$newSheet.Range("H2").Value = 8

$newSheet.Range("B3").Value = "'006569"
$newSheet.Range("C3").Value = "国联安中证医药100指数C"
$newSheet.Range("D3").Value = "'0.32"
$newSheet.Range("E3").Value = "'92.19"
$newSheet.Range("F3").Value = "'1.17"
$newSheet.Range("G3").Value = "'0.0037"
$newSheet.Range("H3").Value = 8

# Restore the original active sheet/tab selection so this edit doesn't leave
# the "2022-Q3" tab selected as a side effect of creating it.
$wb.Worksheets.Item(1).Activate()
